# Updates the Price (D) and Volume(1h) (E) columns of the cryptos list
# to match the latest scrape, per the commit "Updated cryptos list on
# Sat Feb 18 17:43:24 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "24.658.10"
$ws.Cells.Item(2, 5).Value = "  +1.52%  "
$ws.Cells.Item(3, 4).Value = "1.693.84"
$ws.Cells.Item(3, 5).Value = "  +0.25%  "
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  +0.55%  "
$ws.Cells.Item(5, 4).Value = "'316.77"
$ws.Cells.Item(5, 5).Value = "  +2.08%  "
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(7, 4).Value = "'0.3947"
$ws.Cells.Item(7, 5).Value = "  +1.26%  "
$ws.Cells.Item(8, 4).Value = "'0.4051"
$ws.Cells.Item(8, 5).Value = "  +0.56%  "
$ws.Cells.Item(9, 4).Value = "'1.488"
$ws.Cells.Item(9, 5).Value = "  +1.86%  "
$ws.Cells.Item(10, 4).Value = "'1.002"
$ws.Cells.Item(10, 5).Value = "  +0.57%  "
$ws.Cells.Item(11, 4).Value = "'51.88"
$ws.Cells.Item(11, 5).Value = "  -4.69%  "
$ws.Cells.Item(12, 4).Value = "'0.08875"
$ws.Cells.Item(12, 5).Value = "  +1.95%  "
$ws.Cells.Item(13, 4).Value = "'7.163"
$ws.Cells.Item(13, 5).Value = "  -1.41%  "
$ws.Cells.Item(14, 4).Value = "'23.47"
$ws.Cells.Item(14, 5).Value = "  +2.35%  "
$ws.Cells.Item(15, 4).Value = "'8.140"
$ws.Cells.Item(15, 5).Value = "  +9.55%  "
$ws.Cells.Item(16, 4).Value = "'0.00001325"
$ws.Cells.Item(16, 5).Value = "  +0.96%  "
$ws.Cells.Item(17, 4).Value = "1.693.90"
$ws.Cells.Item(17, 5).Value = "  +0.29%  "
$ws.Cells.Item(18, 4).Value = "'99.99"
$ws.Cells.Item(18, 5).Value = "  +0.36%  "
$ws.Cells.Item(19, 4).Value = "'0.07011"
$ws.Cells.Item(19, 5).Value = "  +0.45%  "
$ws.Cells.Item(20, 4).Value = "'19.64"
$ws.Cells.Item(20, 5).Value = "  +2.04%  "
$ws.Cells.Item(21, 4).Value = "'7.014"
$ws.Cells.Item(21, 5).Value = "  +5.20%  "
$ws.Cells.Item(22, 5).Value = "  +0.53%  "
$ws.Cells.Item(23, 4).Value = "'14.32"
$ws.Cells.Item(23, 5).Value = "  +1.87%  "
$ws.Cells.Item(24, 4).Value = "24.647.44"
$ws.Cells.Item(24, 5).Value = "  +1.41%  "
$ws.Cells.Item(25, 4).Value = "'3.200"
$ws.Cells.Item(25, 5).Value = "  +6.32%  "
$ws.Cells.Item(26, 4).Value = "'2.346"
$ws.Cells.Item(27, 4).Value = "'22.69"
$ws.Cells.Item(27, 5).Value = "  +2.88%  "
$ws.Cells.Item(28, 4).Value = "'163.29"
$ws.Cells.Item(28, 5).Value = "  +2.05%  "
$ws.Cells.Item(29, 4).Value = "'136.47"
$ws.Cells.Item(29, 5).Value = "  +3.63%  "
$ws.Cells.Item(30, 4).Value = "'5.161"
$ws.Cells.Item(30, 5).Value = "  +1.43%  "
$ws.Cells.Item(31, 4).Value = "'7.486"
$ws.Cells.Item(31, 5).Value = "  -3.70%  "
$ws.Cells.Item(32, 4).Value = "1.881.10"
$ws.Cells.Item(32, 5).Value = "  +0.33%  "
$ws.Cells.Item(33, 4).Value = "'1.066"
$ws.Cells.Item(33, 5).Value = "  -1.78%  "
$ws.Cells.Item(34, 4).Value = "'0.08584"
$ws.Cells.Item(34, 5).Value = "  -0.58%  "
$ws.Cells.Item(35, 4).Value = "'7.150"
$ws.Cells.Item(35, 5).Value = "  -6.89%  "
$ws.Cells.Item(36, 4).Value = "'11.45"
$ws.Cells.Item(36, 5).Value = "  +1.96%  "
$ws.Cells.Item(37, 4).Value = "'0.2737"
$ws.Cells.Item(37, 5).Value = "  +2.43%  "
$ws.Cells.Item(38, 5).Value = "  -0.83%  "
$ws.Cells.Item(39, 4).Value = "'14.43"
$ws.Cells.Item(39, 5).Value = "  -1.11%  "
$ws.Cells.Item(40, 4).Value = "'0.09167"
$ws.Cells.Item(40, 5).Value = "  +3.44%  "
$ws.Cells.Item(41, 4).Value = "'0.02725"
$ws.Cells.Item(41, 5).Value = "  +0.75%  "
$ws.Cells.Item(42, 4).Value = "'1.476"
$ws.Cells.Item(42, 5).Value = "  +1.58%  "
$ws.Cells.Item(43, 4).Value = "'0.7663"
$ws.Cells.Item(43, 5).Value = "  +1.35%  "
$ws.Cells.Item(44, 4).Value = "'16.05"
$ws.Cells.Item(44, 5).Value = "  +5.14%  "
$ws.Cells.Item(45, 4).Value = "'2.597"
$ws.Cells.Item(45, 5).Value = "  +6.35%  "
$ws.Cells.Item(46, 4).Value = "'0.7183"
$ws.Cells.Item(46, 5).Value = "  +0.96%  "
$ws.Cells.Item(47, 4).Value = "'4.219"
$ws.Cells.Item(47, 5).Value = "  +2.38%  "
$ws.Cells.Item(48, 5).Value = "  +0.57%  "
$ws.Cells.Item(49, 4).Value = "'1.330"
$ws.Cells.Item(49, 5).Value = "  +6.09%  "
$ws.Cells.Item(50, 4).Value = "'140.40"
$ws.Cells.Item(50, 5).Value = "  +0.87%  "
$ws.Cells.Item(51, 4).Value = "'0.07974"
$ws.Cells.Item(51, 5).Value = "  +0.71%  "
